$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.909.74"
$ws.Range("E2").Value = "  +2.04%  "
$ws.Range("D3").Value = "'3.170.15"
$ws.Range("E3").Value = "  +4.13%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'580.31"
$ws.Range("E5").Value = "  +4.11%  "
$ws.Range("D6").Value = "'151.57"
$ws.Range("E6").Value = "  +6.86%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'3.167.28"
$ws.Range("E8").Value = "  +4.15%  "
$ws.Range("E9").Value = "  +2.72%  "
$ws.Range("D10").Value = "'0.162"
$ws.Range("E10").Value = "  +6.21%  "
$ws.Range("D11").Value = "'6.24"
$ws.Range("E11").Value = "  -1.48%  "
$ws.Range("D12").Value = "'0.502"
$ws.Range("E12").Value = "  +3.13%  "
$ws.Range("D13").Value = "'0.0000269"
$ws.Range("E13").Value = "  +17.60%  "
$ws.Range("D14").Value = "'37.73"
$ws.Range("E14").Value = "  +6.07%  "
$ws.Range("D15").Value = "'3.687.70"
$ws.Range("E15").Value = "  +4.19%  "
$ws.Range("D16").Value = "'64.968.06"
$ws.Range("E16").Value = "  +2.01%  "
$ws.Range("E17").Value = "  +5.97%  "
$ws.Range("D18").Value = "'3.166.84"
$ws.Range("E18").Value = "  +4.17%  "
$ws.Range("E19").Value = "  +1.57%  "
$ws.Range("D20").Value = "'512.72"
$ws.Range("E20").Value = "  +8.10%  "
$ws.Range("D21").Value = "'14.85"
$ws.Range("E21").Value = "  +5.74%  "
$ws.Range("D22").Value = "'0.729"
$ws.Range("E22").Value = "  +7.14%  "
$ws.Range("D23").Value = "'15.27"
$ws.Range("E23").Value = "  +4.07%  "
$ws.Range("D24").Value = "'7.81"
$ws.Range("E24").Value = "  +4.27%  "
$ws.Range("D25").Value = "'85.30"
$ws.Range("E25").Value = "  +3.15%  "
$ws.Range("D27").Value = "'9.00"
$ws.Range("E27").Value = "  +10.87%  "
$ws.Range("E28").Value = "  +5.27%  "
$ws.Range("E29").Value = "  +7.39%  "
$ws.Range("D30").Value = "'27.88"
$ws.Range("E30").Value = "  +6.56%  "
$ws.Range("D31").Value = "'2.76"
$ws.Range("E31").Value = "  +13.61%  "
$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("E33").Value = "  +4.70%  "
$ws.Range("D34").Value = "'6.36"
$ws.Range("E34").Value = "  +10.65%  "
$ws.Range("D35").Value = "'6.58"
$ws.Range("E35").Value = "  +6.19%  "
$ws.Range("D36").Value = "'55.69"
$ws.Range("E36").Value = "  +1.98%  "
$ws.Range("D37").Value = "'0.0896"
$ws.Range("E37").Value = "  +10.54%  "
$ws.Range("D38").Value = "'3.15"
$ws.Range("E38").Value = "  +13.21%  "
$ws.Range("D39").Value = "'473.98"
$ws.Range("E39").Value = "  +7.68%  "
$ws.Range("D40").Value = "'0.0419"
$ws.Range("E40").Value = "  +2.65%  "
$ws.Range("D41").Value = "'8.65"
$ws.Range("E41").Value = "  +4.73%  "
$ws.Range("D42").Value = "'3.065.31"
$ws.Range("E43").Value = "  +0.47%  "
$ws.Range("E44").Value = "  +6.15%  "
$ws.Range("E45").Value = "  +7.76%  "
$ws.Range("D46").Value = "'29.15"
$ws.Range("E46").Value = "  +4.89%  "
$ws.Range("E47").Value = "  +19.42%  "
$ws.Range("D48").Value = "'0.999"
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("E49").Value = "  +1.84%  "
$ws.Range("E50").Value = "  +8.43%  "
$ws.Range("D51").Value = "'120.52"
$ws.Range("E51").Value = "  +1.80%  "
